$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algo_programs")

# Row 51 (sheet row 53): the Red-Black-Tree insert program is no longer
# "Testing Left" - it has been verified, so mark it Completed.
$ws.Range("G53").Value = "Completed"

# Row 52 (sheet row 54): new program - deleting a node from a RBT.
$ws.Range("B54").Value = 52
$ws.Range("C54").Value = "delRBTree.c"
$ws.Range("D54").Value = "C"
$ws.Range("E54").Value = "Delete a node from a RBT"
$ws.Range("F54").Value = "Class"
$ws.Range("G54").Value = "Incomplete"

# Row 53 (sheet row 55): new program - AVL insert/delete module.
$ws.Range("B55").Value = 53
$ws.Range("C55").Value = "modAVL.h"
$ws.Range("D55").Value = "C"
$ws.Range("E55").Value = "Code for AVL insert and delete"
$ws.Range("F55").Value = "Class"
$ws.Range("G55").Value = "Completed"

# Row 54 (sheet row 56): new program - testing header for AVL insert/delete.
$ws.Range("B56").Value = 54
$ws.Range("C56").Value = "testing53.c"
$ws.Range("D56").Value = "C"
$ws.Range("E56").Value = "Testing header file for AVL insert and delete"
$ws.Range("F56").Value = "Class"
$ws.Range("G56").Value = "Completed"

# Row 55 (sheet row 57): new program - Red Black Tree insert/delete module.
$ws.Range("B57").Value = 55
$ws.Range("C57").Value = "modRB.h"
$ws.Range("D57").Value = "C"
$ws.Range("E57").Value = "Header file for Red Black Tree insertion and deletion"
$ws.Range("F57").Value = "Class"
$ws.Range("G57").Value = "Completed"

# Row 56 (sheet row 58): new program - testing header for RBT insert/delete.
$ws.Range("B58").Value = 56
$ws.Range("C58").Value = "testing55.c"
$ws.Range("D58").Value = "C"
$ws.Range("E58").Value = "Testing header file for Red Black Tree insertion and deletion"
$ws.Range("F58").Value = "Class"
$ws.Range("G58").Value = "Incomplete"

# Keep the sheet view/selection in sync with where the edits were made.
$ws.Range("G58").Select()
